$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '67.750.13'
$ws.Range("E2").Value = '  +0.59%  '

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.800.57'
$ws.Range("E3").Value = '  +0.43%  '

# Row 4
$ws.Range("E4").Value = '  +0.16%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '596.76'
$ws.Range("E5").Value = '  +0.57%  '

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '167.23'
$ws.Range("E6").Value = '  +0.92%  '

# Row 7
$ws.Range("E7").Value = '  -0.08%  '

# Row 8
$ws.Range("E8").Value = '  +0.41%  '

# Row 9
$ws.Range("E9").Value = '  +1.74%  '

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '6.29'
$ws.Range("E10").Value = '  -1.14%  '

# Row 11
$ws.Range("E11").Value = '  +0.26%  '

# Row 12
$ws.Range("E12").Value = '  -0.25%  '

# Row 13
$ws.Range("E13").Value = '  +0.29%  '

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '4.442.83'
$ws.Range("E14").Value = '  +0.78%  '

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '3.798.23'
$ws.Range("E15").Value = '  +0.89%  '

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '18.56'
$ws.Range("E16").Value = '  +3.41%  '

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '67.800.66'
$ws.Range("E17").Value = '  +0.79%  '

# Row 18
$ws.Range("E18").Value = '  +1.84%  '

# Row 19
$ws.Range("E19").Value = '  +0.51%  '

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '461.32'
$ws.Range("E20").Value = '  +0.86%  '

# Row 21
$ws.Range("E21").Value = '  -2.89%  '

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.699'
$ws.Range("E22").Value = '  -0.08%  '

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.0000152'
$ws.Range("E23").Value = '  +1.58%  '

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '83.27'
$ws.Range("E24").Value = '  -0.05%  '

# Row 25
$ws.Range("E25").Value = '  +2.11%  '

# Row 26
$ws.Range("E26").Value = '  -1.08%  '

# Row 27
$ws.Range("E27").Value = '  +0.01%  '

# Row 28
$ws.Range("E28").Value = '  +0.64%  '

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '3.946.15'
$ws.Range("E29").Value = '  +0.61%  '

# Row 30
$ws.Range("E30").Value = '  -0.14%  '

# Row 31
$ws.Range("B31").Value = 'NEARProtocol'
$ws.Range("C31").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '7.35'
$ws.Range("E31").Value = '  +2.42%  '

# Row 32
$ws.Range("B32").Value = 'ImmutableX'
$ws.Range("C32").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '2.23'
$ws.Range("E32").Value = '  +1.99%  '

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '29.55'
$ws.Range("E33").Value = '  -0.78%  '

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.999'
$ws.Range("E34").Value = '  -0.17%  '

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '9.06'
$ws.Range("E35").Value = '  -1.08%  '

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '3.740.80'
$ws.Range("E36").Value = '  +0.12%  '

# Row 37
$ws.Range("E37").Value = '  +0.60%  '

# Row 38
$ws.Range("E38").Value = '  +2.41%  '

# Row 39
$ws.Range("E39").Value = '  +0.23%  '

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.998'
$ws.Range("E40").Value = '  +0.63%  '

# Row 41
$ws.Range("E41").Value = '  +1.14%  '

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.999'
$ws.Range("E42").Value = '  +0.19%  '

# Row 43
$ws.Range("E43").Value = '  +0.00%  '

# Row 44
$ws.Range("E44").Value = '  +2.50%  '

# Row 45
$ws.Range("E45").Value = '  +1.39%  '

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '42.79'
$ws.Range("E46").Value = '  -2.02%  '

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '8.35'
$ws.Range("E47").Value = '  +0.12%  '

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '27.34'
$ws.Range("E48").Value = '  +7.61%  '

# Row 49
$ws.Range("B49").Value = 'ONDO'
$ws.Range("C49").Value = 'https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.36'
$ws.Range("E49").Value = '  +9.94%  '

# Row 50
$ws.Range("B50").Value = 'Monero'
$ws.Range("C50").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '147.39'
$ws.Range("E50").Value = '  -0.14%  '

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '394.83'
$ws.Range("E51").Value = '  +0.48%  '
